$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 3
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 1

# Row 4
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1

# Row 5
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0

# Row 6
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

# Row 8
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1

# Row 9
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 1

# Row 10
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 1

# Row 11
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 1
